$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set the Price column (D) to Text format first so numeric-looking strings
# (e.g. "42.979.77", "0.0000100", "237.00") are preserved exactly as typed,
# matching the original inline-string formatting used for these cells.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.979.77"
$ws.Range("E2").Value = "  +0.56%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.285.88"
$ws.Range("E3").Value = "  +1.53%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "251.66"
$ws.Range("E5").Value = "  -0.83%  "
$ws.Range("E6").Value = "  +0.86%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "73.27"
$ws.Range("E7").Value = "  +1.88%  "
$ws.Range("E8").Value = "  -0.11%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.645"
$ws.Range("E9").Value = "  +0.52%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.01"
$ws.Range("E10").Value = "  -4.98%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0977"
$ws.Range("E11").Value = "  +1.67%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "59.03"
$ws.Range("E12").Value = "  -0.84%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.41"
$ws.Range("E13").Value = "  +1.16%  "
$ws.Range("E14").Value = "  +1.51%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.628.79"
$ws.Range("E15").Value = "  +1.51%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "15.26"
$ws.Range("E16").Value = "  +3.18%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.868"
$ws.Range("E17").Value = "  -1.90%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.292.54"
$ws.Range("E18").Value = "  +1.81%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "42.861.27"
$ws.Range("E19").Value = "  +0.38%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0000100"
$ws.Range("E20").Value = "  +2.49%  "
$ws.Range("E21").Value = "  +0.68%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "72.49"
$ws.Range("E22").Value = "  -0.58%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "237.00"
$ws.Range("E23").Value = "  +0.82%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.22"
$ws.Range("E24").Value = "  +4.63%  "
$ws.Range("E25").Value = "  -2.27%  "
$ws.Range("E26").Value = "  -0.77%  "
$ws.Range("E27").Value = "  -0.25%  "
$ws.Range("E28").Value = "  -2.22%  "
$ws.Range("E29").Value = "  -0.96%  "
$ws.Range("E30").Value = "  -1.11%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "166.43"
$ws.Range("E31").Value = "  -1.05%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "21.02"
$ws.Range("E32").Value = "  +0.19%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.47"
$ws.Range("E33").Value = "  +4.66%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.127"
$ws.Range("E34").Value = "  -1.80%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0826"
$ws.Range("E35").Value = "  +5.29%  "
$ws.Range("E36").Value = "  +9.35%  "
$ws.Range("E37").Value = "  +1.68%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.60"
$ws.Range("E38").Value = "  +11.16%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.76"
$ws.Range("E39").Value = "  +1.20%  "
$ws.Range("E40").Value = "  -3.31%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "14.20"
$ws.Range("E41").Value = "  +14.35%  "
$ws.Range("E42").Value = "  +1.90%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.91"
$ws.Range("E43").Value = "  +0.00%  "
$ws.Range("E44").Value = "  +7.17%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "9.14"
$ws.Range("E45").Value = "  +2.58%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "61.62"
$ws.Range("E46").Value = "  -4.46%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.87"
$ws.Range("E47").Value = "  -0.54%  "
$ws.Range("E48").Value = "  +1.40%  "
$ws.Range("B49").Value = "Aave"
$ws.Range("C49").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "101.32"
$ws.Range("E49").Value = "  +7.33%  "
$ws.Range("B50").Value = "BinanceUSD"
$ws.Range("C50").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.00"
$ws.Range("E50").Value = "  -0.12%  "
$ws.Range("E51").Value = "  -2.35%  "
